# Weekly update: a new week of "Espinaca" price data (date 44817) is
# prepended to the data table. Two new rows are inserted right after the
# header/existing-row-4 position (i.e. at row 5), shifting the previously
# existing data rows (5-15) down to (7-17), and the two freshly inserted
# rows (5-6) are populated with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 5 (pushes old rows 5..15 down to 7..17,
# carrying their formatting/styles with them).
$ws.Range("A5:R6").EntireRow.Insert()

# New row 5: "Primera" quality entry for the new week (2022-09-13).
$ws.Cells.Item(5,1).Value()  = 7
$ws.Cells.Item(5,2).Value()  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(5,3).Value()  = "Ñuble"
$ws.Cells.Item(5,4).Value()  = 44817
$ws.Cells.Item(5,5).Value()  = 16
$ws.Cells.Item(5,6).Value()  = 100112012
$ws.Cells.Item(5,7).Value()  = "Espinaca"
$ws.Cells.Item(5,8).Value()  = "Sin especificar"
$ws.Cells.Item(5,9).Value()  = "Primera"
$ws.Cells.Item(5,10).Value() = 60
$ws.Cells.Item(5,11).Value() = 7000
$ws.Cells.Item(5,12).Value() = 7000
$ws.Cells.Item(5,13).Value() = 7000
$ws.Cells.Item(5,14).Value() = "`$/cuna 10 kilos"
$ws.Cells.Item(5,15).Value() = "Provincia de Diguillín"
$ws.Cells.Item(5,16).Value() = 700
$ws.Cells.Item(5,17).Value() = 10
$ws.Cells.Item(5,18).Value() = "Hortaliza"

# New row 6: "Segunda" quality entry for the same new week.
$ws.Cells.Item(6,1).Value()  = 7
$ws.Cells.Item(6,2).Value()  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(6,3).Value()  = "Ñuble"
$ws.Cells.Item(6,4).Value()  = 44817
$ws.Cells.Item(6,5).Value()  = 16
$ws.Cells.Item(6,6).Value()  = 100112012
$ws.Cells.Item(6,7).Value()  = "Espinaca"
$ws.Cells.Item(6,8).Value()  = "Sin especificar"
$ws.Cells.Item(6,9).Value()  = "Segunda"
$ws.Cells.Item(6,10).Value() = 60
$ws.Cells.Item(6,11).Value() = 8000
$ws.Cells.Item(6,12).Value() = 8000
$ws.Cells.Item(6,13).Value() = 8000
$ws.Cells.Item(6,14).Value() = "`$/cuna 10 kilos"
$ws.Cells.Item(6,15).Value() = "Provincia de Diguillín"
$ws.Cells.Item(6,16).Value() = 800
$ws.Cells.Item(6,17).Value() = 10
$ws.Cells.Item(6,18).Value() = "Hortaliza"
